$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '36.599.09'
$ws.Range('E2').Value = '  +0.64%  '
$ws.Range('D3').Value = '2.008.19'
$ws.Range('E3').Value = '  -0.24%  '
$ws.Range('E4').Value = '  -0.07%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '247.47'
$ws.Range('E5').Value = '  -1.92%  '
$ws.Range('E6').Value = '  -1.43%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '62.61'
$ws.Range('E7').Value = '  +1.19%  '
$ws.Range('E8').Value = '  -0.01%  '
$ws.Range('E9').Value = '  +3.91%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '57.01'
$ws.Range('E10').Value = '  -1.98%  '
$ws.Range('E11').Value = '  +4.87%  '
$ws.Range('E12').Value = '  -0.08%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '0.887'
$ws.Range('E13').Value = '  -2.21%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '22.64'
$ws.Range('E14').Value = '  +10.45%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '14.13'
$ws.Range('E15').Value = '  -5.46%  '
$ws.Range('D16').Value = '2.299.33'
$ws.Range('E16').Value = '  -0.36%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '5.51'
$ws.Range('E17').Value = '  +0.58%  '
$ws.Range('D18').Value = '2.009.06'
$ws.Range('E18').Value = '  -0.18%  '
$ws.Range('D19').Value = '36.513.56'
$ws.Range('E19').Value = '  +0.46%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '71.97'
$ws.Range('E20').Value = '  -0.08%  '
$ws.Range('D21').Value = '0.0₃0872'
$ws.Range('E21').Value = '  +0.75%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '5.32'
$ws.Range('E22').Value = '  +0.33%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '238.96'
$ws.Range('E23').Value = '  +1.86%  '
$ws.Range('E24').Value = '  +0.05%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '2.51'
$ws.Range('E25').Value = '  -7.55%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '2.33'
$ws.Range('E26').Value = '  +0.67%  '
$ws.Range('E27').Value = '  +2.92%  '
$ws.Range('B28').Value = 'Monero'
$ws.Range('C28').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '159.79'
$ws.Range('E28').Value = '  -2.33%  '
$ws.Range('B29').Value = 'Kaspa'
$ws.Range('C29').Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '0.136'
$ws.Range('E29').Value = '  +22.58%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '20.20'
$ws.Range('E30').Value = '  +2.82%  '
$ws.Range('E31').Value = '  +0.74%  '
$ws.Range('B32').Value = 'Filecoin'
$ws.Range('C32').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '5.01'
$ws.Range('E32').Value = '  -2.31%  '
$ws.Range('B33').Value = 'ImmutableX'
$ws.Range('C33').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '1.18'
$ws.Range('E33').Value = '  -0.68%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '0.0625'
$ws.Range('E34').Value = '  +2.49%  '
$ws.Range('E35').Value = '  -2.36%  '
$ws.Range('E36').Value = '  +10.01%  '
$ws.Range('E37').Value = '  -3.61%  '
$ws.Range('E38').Value = '  +0.03%  '
$ws.Range('E39').Value = '  +0.91%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '3.21'
$ws.Range('E40').Value = '  +19.49%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '1.28'
$ws.Range('E41').Value = '  +3.82%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.101'
$ws.Range('E42').Value = '  -3.90%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '2.91'
$ws.Range('E43').Value = '  -1.38%  '
$ws.Range('E44').Value = '  -0.64%  '
$ws.Range('E45').Value = '  -0.93%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '16.78'
$ws.Range('E46').Value = '  -1.73%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '93.67'
$ws.Range('E47').Value = '  -1.60%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '7.67'
$ws.Range('E48').Value = '  -5.44%  '
$ws.Range('D49').Value = '1.360.83'
$ws.Range('E49').Value = '  -6.44%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '2.89'
$ws.Range('E50').Value = '  -1.63%  '
$ws.Range('D51').Value = '2.195.28'
$ws.Range('E51').Value = '  -0.10%  '
